$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update status-like values in the "Reviewer ETA" column (E) ---
# Row 5  : OData Core - Notational Conventions  -> was "ongoing", now "In Progress"
# Row 7  : OData Core - Extensibility           -> was a date, now "Done"
# Row 16 : OData Core - Appendix - 2 initial sections -> was a date, now "Done"
$ws.Range("E5").Value = "In Progress"
$ws.Range("E7").Value = "Done"
$ws.Range("E16").Value = "Done"

# --- New conditional formatting on column E (Reviewer ETA) ---

# E2:E25 -> highlight "DONE" with the standard green "Good" style
$rngE2 = $ws.Range("E2:E25")
$cfDoneGood = $rngE2.FormatConditions.Add(1, 3, '"DONE"')
$cfDoneGood.Font.Color = 24832
$cfDoneGood.Interior.Color = 13561798

# E1:E1048576 -> highlight "In Progress" (amber) and "DONE" (red)
$rngEAll = $ws.Range("E1:E1048576")
$cfInProgress = $rngEAll.FormatConditions.Add(1, 3, '"In Progress"')
$cfInProgress.Interior.Color = 49407

$cfDoneRed = $rngEAll.FormatConditions.Add(1, 3, '"DONE"')
$cfDoneRed.Interior.Color = 255

# --- Renumber conditional formatting priorities so the 3 new rules sit
#     at the top (priority 1-3) and every pre-existing rule is pushed
#     down by 3, matching Excel's normal "newest rule wins" behaviour ---
$ws.Range("G1:G16").FormatConditions.Item(1).Priority = 9
$ws.Range("G1:G16").FormatConditions.Item(2).Priority = 10
$ws.Range("G1:G16").FormatConditions.Item(3).Priority = 11
$ws.Range("C1:C16").FormatConditions.Item(1).Priority = 8
$ws.Range("G17").FormatConditions.Item(1).Priority = 5
$ws.Range("G17").FormatConditions.Item(2).Priority = 6
$ws.Range("G17").FormatConditions.Item(3).Priority = 7
$ws.Range("C17").FormatConditions.Item(1).Priority = 4

$cfDoneGood.Priority = 3
$cfInProgress.Priority = 2
$cfDoneRed.Priority = 1

# --- Final selection left on F24 ---
$ws.Range("F24").Select() | Out-Null
